# Updated cryptos list on Tue Aug  8 02:20:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric/percentage columns for rows 2-47 (and 28, 51 E-only)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.150.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6182"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07343"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2872"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.05"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.944"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6630"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.03"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008938"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.825"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.120.32"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.070.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.17"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.312"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1418"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.480"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.481"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05880"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.080"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.059"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.207"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.862"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7293"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.610"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.847"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.213.33"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01747"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9169"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.262"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.973.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5082"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.39%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.01%  "

# Rows 48-50 coin order changed (EnergySwap, TheSandbox, BabyDogeCoin)
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.105"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4008"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000116"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.90%  "
